# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 502.25
$ws.Range("I2").Value = 579.6667
$ws.Range("K2").Value = 579.6667
$ws.Range("M2").Value = -466.6667
$ws.Range("H17").Value = 2590.8
$ws.Range("I17").Value = 3109.25
$ws.Range("J17").Value = 1998.2858
$ws.Range("K17").Value = 9327.75
$ws.Range("L17").Value = 5994.857400000001
$ws.Range("M17").Value = -9159.75
$ws.Range("N17").Value = -6330.857400000001
$ws.Range("H107").Value = 379.6
$ws.Range("I107").Value = 174.75
$ws.Range("K107").Value = 174.75
$ws.Range("M107").Value = 1745.25
$ws.Range("H111").Value = 2828.1667
$ws.Range("I111").Value = 2828.1667
$ws.Range("K111").Value = 8484.500100000001
$ws.Range("M111").Value = -5417.500100000001
$ws.Range("H112").Value = 2397.5
$ws.Range("I112").Value = 932.6667
$ws.Range("J112").Value = 2690.4666
$ws.Range("K112").Value = 2798.0001
$ws.Range("L112").Value = 8071.399800000001
$ws.Range("M112").Value = -1690.0001
$ws.Range("N112").Value = -10287.3998
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1897.75
$ws.Range("I61").Value = 1897.75
$ws.Range("K61").Value = 1897.75
$ws.Range("M61").Value = -1685.75
$ws.Range("H110").Value = 10333
$ws.Range("I110").Value = 15199.5
$ws.Range("J110").Value = 600
$ws.Range("K110").Value = 15199.5
$ws.Range("L110").Value = 600
$ws.Range("M110").Value = -13154.5
$ws.Range("N110").Value = -4690
$ws.Range("H132").Value = 2579.2856
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H136").Value = 1897.75
$ws.Range("I136").Value = 1897.75
$ws.Range("K136").Value = 5693.25
$ws.Range("M136").Value = -3143.25
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 880.4286
$ws.Range("I64").Value = 790.4286
$ws.Range("K64").Value = 790.4286
$ws.Range("M64").Value = -565.4286
$ws.Range("H67").Value = 880.4286
$ws.Range("I67").Value = 790.4286
$ws.Range("K67").Value = 790.4286
$ws.Range("M67").Value = -10.42859999999996
$ws.Range("H134").Value = 2440.5293
$ws.Range("I134").Value = 2177.8572
$ws.Range("K134").Value = 6533.571599999999
$ws.Range("M134").Value = -3998.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4598
$ws.Range("I31").Value = 1956
$ws.Range("J31").Value = 7900.5
$ws.Range("K31").Value = 1956
$ws.Range("L31").Value = 7900.5
$ws.Range("M31").Value = -1661
$ws.Range("N31").Value = -8490.5
$ws.Range("H34").Value = 4598
$ws.Range("I34").Value = 1956
$ws.Range("J34").Value = 7900.5
$ws.Range("K34").Value = 1956
$ws.Range("L34").Value = 7900.5
$ws.Range("M34").Value = -1754
$ws.Range("N34").Value = -8304.5
$ws.Range("H58").Value = 3673.75
$ws.Range("I58").Value = 1812.6666
$ws.Range("J58").Value = 9257
$ws.Range("K58").Value = 1812.6666
$ws.Range("L58").Value = 9257
$ws.Range("M58").Value = -1609.6666
$ws.Range("N58").Value = -9663
$ws.Range("H99").Value = 13334
$ws.Range("I99").Value = 10392.333
$ws.Range("K99").Value = 10392.333
$ws.Range("M99").Value = -8894.333000000001
$ws.Range("H107").Value = 1080.4546
$ws.Range("I107").Value = 854.6667
$ws.Range("J107").Value = 1351.4
$ws.Range("K107").Value = 854.6667
$ws.Range("L107").Value = 1351.4
$ws.Range("M107").Value = 1065.3333
$ws.Range("N107").Value = -5191.4
$ws.Range("H122").Value = 3332.4167
$ws.Range("I122").Value = 3451.5217
$ws.Range("J122").Value = 3121.6924
$ws.Range("K122").Value = 10354.5651
$ws.Range("L122").Value = 9365.0772
$ws.Range("M122").Value = -7904.5651
$ws.Range("N122").Value = -14265.0772
$ws.Range("H126").Value = 13334
$ws.Range("I126").Value = 10392.333
$ws.Range("K126").Value = 31176.999
$ws.Range("M126").Value = -28706.999
$ws.Range("H132").Value = 2646.3333
$ws.Range("I132").Value = 1565.9166
$ws.Range("J132").Value = 6968
$ws.Range("K132").Value = 4697.7498
$ws.Range("L132").Value = 20904
$ws.Range("M132").Value = -2167.7498
$ws.Range("N132").Value = -25964
$ws.Range("H134").Value = 3796.1333
$ws.Range("I134").Value = 3047.5
$ws.Range("K134").Value = 9142.5
$ws.Range("M134").Value = -6607.5
$ws.Range("H136").Value = 3673.75
$ws.Range("I136").Value = 1812.6666
$ws.Range("J136").Value = 9257
$ws.Range("K136").Value = 5437.9998
$ws.Range("L136").Value = 27771
$ws.Range("M136").Value = -2887.9998
$ws.Range("N136").Value = -32871
$ws.Range("H141").Value = 20546.562
$ws.Range("J141").Value = 20546.562
$ws.Range("L141").Value = 20546.562
$ws.Range("N141").Value = -30906.562

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66746
$ws.Range("J2").Value = 171.2
$ws.Range("L2").Value = 1027.2
$ws.Range("N2").Value = -1253.2
$ws.Range("H39").Value = 756.6667
$ws.Range("I39").Value = 756.6667
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2270.0001
$ws.Range("L39").Value = 0
$ws.Range("H131").Value = 1732.9183
$ws.Range("I131").Value = 865
$ws.Range("J131").Value = 1877.5714
$ws.Range("K131").Value = 2595
$ws.Range("L131").Value = 5632.7142
$ws.Range("M131").Value = 2445
$ws.Range("N131").Value = -15712.7142
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 11999.667
$ws.Range("J36").Value = 10999.5
$ws.Range("L36").Value = 10999.5
$ws.Range("N36").Value = -11969.5
$ws.Range("H102").Value = 2044.2
$ws.Range("I102").Value = 535.3333
$ws.Range("J102").Value = 4307.5
$ws.Range("K102").Value = 535.3333
$ws.Range("L102").Value = 4307.5
$ws.Range("M102").Value = 1086.6667
$ws.Range("N102").Value = -7551.5
$ws.Range("H122").Value = 86749.914
$ws.Range("I122").Value = 2828.111
$ws.Range("K122").Value = 8484.332999999999
$ws.Range("M122").Value = -6034.332999999999
$ws.Range("H126").Value = 3808.5454
$ws.Range("J126").Value = 4399.3335
$ws.Range("L126").Value = 13198.0005
$ws.Range("N126").Value = -18138.0005
$ws.Range("H132").Value = 1979.9333
$ws.Range("I132").Value = 1542.7693
$ws.Range("J132").Value = 4821.5
$ws.Range("K132").Value = 4628.3079
$ws.Range("L132").Value = 14464.5
$ws.Range("M132").Value = -2098.3079
$ws.Range("N132").Value = -19524.5
$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 240000
$ws.Range("N134").Value = -245070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2295.9092
$ws.Range("I7").Value = 2139.7778
$ws.Range("K7").Value = 2139.7778
$ws.Range("M7").Value = -2027.7778
$ws.Range("H46").Value = 3717.3125
$ws.Range("I46").Value = 1745
$ws.Range("K46").Value = 1745
$ws.Range("M46").Value = -1557
$ws.Range("H122").Value = 11282.071
$ws.Range("I122").Value = 11079.333
$ws.Range("K122").Value = 33237.999
$ws.Range("M122").Value = -30787.999
$ws.Range("H126").Value = 2295.9092
$ws.Range("I126").Value = 2139.7778
$ws.Range("K126").Value = 6419.3334
$ws.Range("M126").Value = -3949.3334
$ws.Range("H132").Value = 3758.375
$ws.Range("I132").Value = 3557.1904
$ws.Range("K132").Value = 10671.5712
$ws.Range("M132").Value = -8141.5712
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 860.3
$ws.Range("I81").Value = 917
$ws.Range("K81").Value = 1834
$ws.Range("M81").Value = -773
$ws.Range("H84").Value = 860.3
$ws.Range("I84").Value = 917
$ws.Range("K84").Value = 9170
$ws.Range("M84").Value = -3866
$ws.Range("H136").Value = 2900.2856
$ws.Range("I136").Value = 1089.1538
$ws.Range("K136").Value = 3267.4614
$ws.Range("M136").Value = -717.4614000000001
